$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as literal
# text (so values like "1.000" or "29.406.85" are preserved exactly instead
# of being reinterpreted as numbers), then restore the default "Normal" style
# so the cell format matches the source workbook (General, no quote-prefix).
function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# --- Update price (D) and volume/1h change (E) for each coin row ---

Set-TextValue "D2" "29.406.85"
Set-TextValue "E2" "  -0.30%  "

Set-TextValue "D3" "1.849.55"
Set-TextValue "E3" "  -0.16%  "

Set-TextValue "D4" "0.9990"
Set-TextValue "E4" "  -0.05%  "

Set-TextValue "D5" "240.67"

Set-TextValue "D6" "0.6312"
Set-TextValue "E6" "  -0.22%  "

Set-TextValue "D7" "1.000"
Set-TextValue "E7" "  +0.09%  "

Set-TextValue "D8" "0.07568"
Set-TextValue "E8" "  +0.44%  "

Set-TextValue "D9" "0.2959"
Set-TextValue "E9" "  -0.76%  "

Set-TextValue "D10" "24.41"
Set-TextValue "E10" "  -0.39%  "

Set-TextValue "D11" "0.07703"
Set-TextValue "E11" "  -0.27%  "

Set-TextValue "D12" "1.856.86"
Set-TextValue "E12" "  -0.34%  "

Set-TextValue "D13" "5.000"
Set-TextValue "E13" "  -0.29%  "

Set-TextValue "D14" "0.6853"
Set-TextValue "E14" "  -1.40%  "

Set-TextValue "D15" "0.00001002"
Set-TextValue "E15" "  +1.99%  "

Set-TextValue "D16" "83.13"
Set-TextValue "E16" "  -0.46%  "

Set-TextValue "D17" "2.101.56"
Set-TextValue "E17" "  -1.65%  "

Set-TextValue "D18" "6.152"
Set-TextValue "E18" "  -1.97%  "

Set-TextValue "D19" "29.431.36"
Set-TextValue "E19" "  -0.43%  "

Set-TextValue "D20" "227.97"
Set-TextValue "E20" "  -2.44%  "

Set-TextValue "D21" "12.47"
Set-TextValue "E21" "  -0.46%  "

Set-TextValue "D22" "1.000"
Set-TextValue "E22" "  +0.10%  "

Set-TextValue "D23" "7.571"
Set-TextValue "E23" "  -1.30%  "

Set-TextValue "D24" "1.001"
Set-TextValue "E24" "  +0.01%  "

Set-TextValue "D25" "157.24"
Set-TextValue "E25" "  +1.20%  "

Set-TextValue "D26" "0.1396"
Set-TextValue "E26" "  +0.00%  "

Set-TextValue "D27" "8.387"
Set-TextValue "E27" "  -0.96%  "

Set-TextValue "D28" "17.70"
Set-TextValue "E28" "  -0.16%  "

Set-TextValue "D29" "1.470"
Set-TextValue "E29" "  -0.40%  "

Set-TextValue "D30" "0.05711"
Set-TextValue "E30" "  -3.62%  "

Set-TextValue "D31" "1.258"
Set-TextValue "E31" "  +0.44%  "

Set-TextValue "D32" "4.126"
Set-TextValue "E32" "  -0.05%  "

Set-TextValue "D33" "4.025"
Set-TextValue "E33" "  -0.08%  "

Set-TextValue "D34" "1.847"
Set-TextValue "E34" "  -2.56%  "

Set-TextValue "D35" "1.156"
Set-TextValue "E35" "  -0.97%  "

Set-TextValue "D36" "0.7181"
Set-TextValue "E36" "  -0.74%  "

Set-TextValue "D37" "2.592"
Set-TextValue "E37" "  +0.40%  "

Set-TextValue "D38" "1.252.13"
Set-TextValue "E38" "  +1.16%  "

Set-TextValue "D39" "0.01809"
Set-TextValue "E39" "  +0.82%  "

Set-TextValue "D40" "2.782"
Set-TextValue "E40" "  -0.54%  "

Set-TextValue "D41" "0.9062"
Set-TextValue "E41" "  -0.27%  "

Set-TextValue "D42" "6.183"
Set-TextValue "E42" "  +1.30%  "

Set-TextValue "E43" "  +0.19%  "

Set-TextValue "D44" "101.43"
Set-TextValue "E44" "  -0.02%  "

Set-TextValue "D45" "66.17"
Set-TextValue "E45" "  -1.96%  "

Set-TextValue "D49" "9.108"
Set-TextValue "E49" "  -0.71%  "

Set-TextValue "D50" "1.682"
Set-TextValue "E50" "  -1.60%  "

Set-TextValue "E51" "  +0.06%  "

# --- Rows 46-48 were re-ranked: Aptos (was row 46) dropped below BabyDogeCoin
# (was row 48), and TheSandbox (was row 47) dropped to row 48. Update Coin,
# Link, Price and Volume(1h) for the new ranking order. ---
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D46" "0.00000000120"
Set-TextValue "E46" "  +0.68%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D47" "7.104"
Set-TextValue "E47" "  -3.29%  "

$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D48" "0.4041"
Set-TextValue "E48" "  -0.17%  "
